# Apply edit: add columns I (I0) and J (IF) with per-row integer values,
# matching the commit "I0 and IF added".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for the two new columns, matching style of existing headers (s="1" -> bold/centered/bordered)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Data rows 2-92: row index, I value, J value
$data = @(
    @(2,7,8),
    @(3,7,7),
    @(4,8,8),
    @(5,7,8),
    @(6,7,8),
    @(7,8,8),
    @(8,7,7),
    @(9,6,7),
    @(10,7,7),
    @(11,8,8),
    @(12,7,7),
    @(13,7,8),
    @(14,8,8),
    @(15,5,5),
    @(16,8,8),
    @(17,6,7),
    @(18,8,8),
    @(19,7,7),
    @(20,7,7),
    @(21,10,10),
    @(22,7,7),
    @(23,8,8),
    @(24,7,7),
    @(25,7,7),
    @(26,8,8),
    @(27,7,8),
    @(28,7,8),
    @(29,6,6),
    @(30,8,8),
    @(31,7,7),
    @(32,7,8),
    @(33,6,6),
    @(34,7,7),
    @(35,8,8),
    @(36,7,7),
    @(37,8,8),
    @(38,7,7),
    @(39,8,8),
    @(40,6,6),
    @(41,7,7),
    @(42,6,7),
    @(43,6,6),
    @(44,6,7),
    @(45,8,8),
    @(46,7,7),
    @(47,9,9),
    @(48,5,5),
    @(49,9,9),
    @(50,8,8),
    @(51,9,9),
    @(52,9,9),
    @(53,3,4),
    @(54,6,7),
    @(55,7,7),
    @(56,6,6),
    @(57,8,8),
    @(58,7,7),
    @(59,6,6),
    @(60,7,7),
    @(61,9,9),
    @(62,7,7),
    @(63,8,8),
    @(64,8,8),
    @(65,8,8),
    @(66,8,8),
    @(67,7,7),
    @(68,7,7),
    @(69,9,9),
    @(70,7,7),
    @(71,8,8),
    @(72,8,8),
    @(73,6,7),
    @(74,7,8),
    @(75,8,8),
    @(76,9,9),
    @(77,8,8),
    @(78,10,10),
    @(79,8,8),
    @(80,7,7),
    @(81,7,7),
    @(82,6,6),
    @(83,6,6),
    @(84,4,4),
    @(85,7,7),
    @(86,9,9),
    @(87,9,9),
    @(88,8,8),
    @(89,6,7),
    @(90,7,7),
    @(91,3,3),
    @(92,6,6)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 9).Value = $row[1]
    $ws.Cells.Item($r, 10).Value = $row[2]
}
